$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed cryptos list: rows 37/38 (ImmutableX <-> FirstDigitalUSD) were swapped,
# and every row got updated Price (D) / Volume(1h) (E) values from the latest fetch.
# Price values look numeric, so plain assignment would make Excel coerce them into
# real numbers (e.g. "1.00" -> 1, dropping the trailing zero). To keep them as the
# literal text the source data uses, force the cell to Text format first, then
# restore the original (default) cell style so no extra formatting is introduced.

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "58.768.03"
$ws.Range("E2").Value = "  -0.30%  "
Set-TextValue "D3" "2.309.02"
$ws.Range("E3").Value = "  -0.34%  "
$ws.Range("E4").Value = "  +0.02%  "
Set-TextValue "D5" "538.58"
$ws.Range("E5").Value = "  -1.91%  "
Set-TextValue "D6" "132.52"
$ws.Range("E6").Value = "  +0.67%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +2.40%  "
Set-TextValue "D9" "2.306.98"
$ws.Range("E9").Value = "  -0.40%  "
$ws.Range("E10").Value = "  -1.65%  "
Set-TextValue "D11" "5.49"
$ws.Range("E11").Value = "  -1.40%  "
$ws.Range("E12").Value = "  +1.00%  "
$ws.Range("E13").Value = "  -0.52%  "
Set-TextValue "D14" "23.76"
$ws.Range("E14").Value = "  -1.12%  "
Set-TextValue "D15" "2.721.77"
$ws.Range("E15").Value = "  -0.35%  "
Set-TextValue "D16" "58.645.05"
$ws.Range("E16").Value = "  -0.38%  "
$ws.Range("E17").Value = "  -0.41%  "
Set-TextValue "D18" "2.331.70"
$ws.Range("E18").Value = "  -0.23%  "
$ws.Range("E19").Value = "  -1.04%  "
Set-TextValue "D20" "4.19"
$ws.Range("E20").Value = "  -3.30%  "
Set-TextValue "D21" "314.06"
$ws.Range("E21").Value = "  -0.66%  "
Set-TextValue "D22" "6.64"
$ws.Range("E22").Value = "  +1.95%  "
Set-TextValue "D24" "62.91"
$ws.Range("E24").Value = "  -1.07%  "
$ws.Range("E25").Value = "  +0.48%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("E27").Value = "  -1.92%  "
$ws.Range("E28").Value = "  -1.87%  "
Set-TextValue "D29" "171.42"
$ws.Range("E29").Value = "  +1.07%  "
$ws.Range("E30").Value = "  -2.27%  "
Set-TextValue "D31" "0.0₃0735"
$ws.Range("E31").Value = "  +0.49%  "
Set-TextValue "D32" "1.15"
$ws.Range("E32").Value = "  +2.44%  "
$ws.Range("E33").Value = "  +0.95%  "
$ws.Range("E34").Value = "  +0.61%  "
$ws.Range("E35").Value = "  +0.67%  "
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D37" "1.29"
$ws.Range("E37").Value = "  +2.21%  "
$ws.Range("B38").Value = "FirstDigitalUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D38" "1.00"
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("E39").Value = "  +1.34%  "
$ws.Range("E40").Value = "  -0.06%  "
Set-TextValue "D41" "295.10"
$ws.Range("E41").Value = "  -2.60%  "
Set-TextValue "D42" "141.29"
$ws.Range("E42").Value = "  -0.25%  "
$ws.Range("E44").Value = "  +0.90%  "
Set-TextValue "D45" "0.0497"
$ws.Range("E45").Value = "  -1.33%  "
Set-TextValue "D46" "0.556"
$ws.Range("E46").Value = "  -0.71%  "
Set-TextValue "D47" "18.39"
$ws.Range("E47").Value = "  -1.87%  "
$ws.Range("E48").Value = "  -2.20%  "
Set-TextValue "D49" "10.99"
$ws.Range("E49").Value = "  -0.38%  "
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("E51").Value = "  +0.63%  "
